# updateExcel: fill in QA results gathered by updateToDb func
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "No"
$ws.Range("E2").Value = "No"

# Row 3
$ws.Range("B3").Value = "NA"
$ws.Range("C3").Value = "NA"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("H3").Value = "NA"
$ws.Range("I3").Value = "NA"
$ws.Range("J3").Value = "NA"

# Row 4
$ws.Range("B4").Value = "Yes"
$ws.Range("K4").Value = "P4:INC000005135152: Wrong schedule on 'True Crime' (154)[Last Update: Raised a new GN case : 01575811 to address title mismatch.]"

# Row 5
$ws.Range("B5").Value = "No"

# Row 6
$ws.Range("B6").Value = "No"
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "NA"

# Row 7
$ws.Range("B7").Value = "No"
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "NA"

# Row 8
$ws.Range("B8").Value = "No"
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "NA"

# Row 9
$ws.Range("B9").Value = "No"

# Row 10
$ws.Range("B10").Value = "No"

# Row 11
$ws.Range("B11").Value = "No"
$ws.Range("C11").Value = "Yes"
$ws.Range("J11").Value = "Yes"
$ws.Range("K11").Value = "P4: INC000005135676: 01573557|[IE]DTV|14:00|Poster Missing in On Demand | Prime Boxset | Kids+Movies [Last Update: VNOC seeking for an update from GN for the remaining posters to be published]"

# Row 12
$ws.Range("B12").Value = "No"

# Row 13
$ws.Range("B13").Value = "No"

# Row 14
$ws.Range("B14").Value = "No"

# Row 15
$ws.Range("B15").Value = "No"
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("G15").Value = "NA"

# Row 16
$ws.Range("B16").Value = "No"
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("G16").Value = "NA"

# Row 17
$ws.Range("B17").Value = "No"
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("G17").Value = "NA"

# Row 18
$ws.Range("B18").Value = "No"

# Row 19
$ws.Range("B19").Value = "No"
$ws.Range("C19").Value = "NA"
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("H19").Value = "NA"
$ws.Range("I19").Value = "NA"
